$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update column C (the "Förändrad" date column) for rows 2 through 27
# from Excel date serial 45259 (2023-11-29) to 45260 (2023-11-30).
for ($row = 2; $row -le 27; $row++) {
    $ws.Cells.Item($row, 3).Value = 45260
}
